# ShowBookings View und Methode im BookingController erstellt
# Mark the "ShowBookings View erstellen" (row 57) and
# "ShowBookings Methode im BookingController erstellen" (row 58) tasks as done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("D56").Copy()
$ws.Range("D57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D58").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C57").Value = "done"
$ws.Range("D57").Value = 43567

$ws.Range("C58").Value = "done"
$ws.Range("D58").Value = 43567

$ws.Range("D58").Select()
